# Backlog_3.xlsx edit: replace the "Semana 03" week label (text) in column C
# ("Semana") of both the SPN and ITI sheets with the plain numeric week
# value 3, then update the active sheet / selection to reflect the user's
# navigation (SPN tab selected, column C highlighted on both sheets).

$wb = $excel.ActiveWorkbook
$wsSPN = $wb.Worksheets.Item("SPN")
$wsITI = $wb.Worksheets.Item("ITI")

# --- ITI sheet: column C already uses the numeric-style formatting (style
# index 16), so we only need to replace the text value with the number 3.
$wsITI.Range("C2:C22").Value = 3

# --- SPN sheet: column C currently uses the "text" style (style index 1).
# Copy the number formatting/style already used on ITI!C2 (style index 16)
# onto SPN!C2:C27 before writing the numeric value, so the cell style
# matches exactly (same font/alignment as the rest of the workbook's week
# numbers) instead of creating a brand-new style entry.
$wsITI.Range("C2").Copy()
$dstSPN = $wsSPN.Range("C2:C27")
$dstSPN.PasteSpecial(-4122)   # xlPasteFormats
$dstSPN.Value = 3

$excel.CutCopyMode = 0

# --- Update sheet views/selections to match the saved workbook state:
# ITI: select C2:C22 (while it's the active sheet) ...
$wsITI.Activate()
$wsITI.Range("C2:C22").Select()

# ... then activate SPN and select C2:C27, leaving SPN as the active
# (tab-selected) sheet, scrolled back to show column A again.
$wsSPN.Activate()
$wsSPN.Range("C2:C27").Select()
